# Append a new row 96 to each of the 4 worksheets, mirroring the
# existing data pattern (columns A-I), matching the committed diff.

$wb = $excel.ActiveWorkbook

$rows = @(
    @{
        Sheet = 1
        A = 45882.43866898148
        B = "0x01,0x7c"
        C = "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,"
        D = "0x01,0x24"
        E = "0x14"
        F = 380
        G = "7.598631275147109e+23"
        H = 292
        I = 14
    },
    @{
        Sheet = 2
        A = 45882.43866898148
        B = "0x01,0x7c"
        C = "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,"
        D = "0x01,0x2C"
        E = "0xe"
        F = 380
        G = "5.68432987514711e+23"
        H = 300
        I = 14
    },
    @{
        Sheet = 3
        A = 45882.43866898148
        B = "0x00,0x82"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
        D = "0x00,0x75"
        E = "0x7"
        F = 130
        G = "5.68631262647114e+23"
        H = 117
        I = 7
    },
    @{
        Sheet = 4
        A = 45882.43866898148
        B = "0x00,0x82"
        C = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
        D = "0x00,0x74"
        E = "0x3"
        F = 130
        G = "9.85046333984776e+23"
        H = 116
        I = 3
    }
)

foreach ($r in $rows) {
    $ws = $wb.Worksheets.Item($r.Sheet)

    $ws.Range("A96").Value = $r.A
    $ws.Range("A96").NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Range("B96").Value = $r.B
    $ws.Range("C96").Value = $r.C
    $ws.Range("D96").Value = $r.D
    $ws.Range("E96").Value = $r.E

    $ws.Range("F96").Value = $r.F
    $ws.Range("G96").Value = [double]$r.G
    $ws.Range("H96").Value = $r.H
    $ws.Range("I96").Value = $r.I
}
